# Sleep Diary workbook update:
#  - Fill in Week 2 table (rows 25-38) with data for all 7 days
#  - Add a new Week 3 table (rows 40-57), copied from the Week 2 block's
#    layout/formatting, with data filled in for the first 3 days only
#  - Minor style tweaks (bold cell for "无" row, special border cell)
#  - Update the active selection / page setup to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Duplicate the Week-2 block (rows 21-38) down to rows 40-57 to create
#    the Week-3 table, before we fill in any new data, so the copy is
#    still blank in the data cells (matches the "3 of 7 days filled"
#    target state once we only populate columns B:D below).
# ---------------------------------------------------------------------
$srcRows = $ws.Range("A21:H38")
$dstRows = $ws.Range("A40:H57")
$srcRows.Copy($dstRows)

# Row heights don't come along with the COM Copy of row ranges in this
# host, so re-apply them explicitly from the source rows.
for ($r = 40; $r -le 57; $r++) {
    $srcRowIndex = $r - 19
    $h = $ws.Rows($srcRowIndex).RowHeight()
    $ws.Rows($r).RowHeight = $h
}

Write-Host "Week 3 block copied"

# ---------------------------------------------------------------------
# 2) Fill in the Week-2 table (rows 25-38) for all 7 days (columns B:H).
# ---------------------------------------------------------------------
function Fill-Row($rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, 2 + $i).Value = $values[$i]
    }
}

Fill-Row 25 @("7：30","7：14","7：27","7：23","6：46","8：12","6：04")
Fill-Row 26 @("7：30","7：30","7：40","7：40","7：20","9：30","10：00")
Fill-Row 27 @("22：50","22：00","22：00","23：00","23：00","23：20","23：30")
Fill-Row 28 @("11：00","24：00","23：40","23：20","23：30","24：20","24：00")
Fill-Row 29 @(20,10,5,10,20,20,40)
Fill-Row 30 @(2,2,3,2,3,0,1)
Fill-Row 31 @(14,5,15,10,10,0,30)
Fill-Row 32 @(480,420,450,480,400,430,280)
Fill-Row 33 @("无","无","无","无","无","无","无")
Fill-Row 34 @(30,60,40,15,10,40,40)
Fill-Row 35 @(2,1,3,2,2,4,3)
Fill-Row 36 @(3,4,3,4,2,3,4)
Fill-Row 37 @(3,2,3,3,4,4,4)
Fill-Row 38 @("无","无","无","无","无","无","有 30")

Write-Host "Week 2 data filled"

# ---------------------------------------------------------------------
# 2b) Small style touch-ups that were left behind in the Week-2 table
#     while the data was being typed in:
#      - H27 picked up a "boxed" left/right border (no top/bottom)
#      - E33:H33 ("无" answers) picked up a bold 微软雅黑 font
# ---------------------------------------------------------------------
$h27 = $ws.Range("H27")
$h27.Borders.Item(7).LineStyle = 1        # xlEdgeLeft
$h27.Borders.Item(7).Weight = 2           # xlThin
$h27.Borders.Item(7).ColorIndex = -4105   # xlColorIndexAutomatic
$h27.Borders.Item(10).LineStyle = 1       # xlEdgeRight
$h27.Borders.Item(10).Weight = 2          # xlThin
$h27.Borders.Item(10).ColorIndex = -4105  # xlColorIndexAutomatic
$h27.Borders.Item(8).LineStyle = -4142    # xlEdgeTop -> none
$h27.Borders.Item(9).LineStyle = -4142    # xlEdgeBottom -> none

$boldRange = $ws.Range("E33:H33")
$boldRange.Font.Name = "微软雅黑"
$boldRange.Font.Bold = $true
$boldRange.Font.Size = 11

Write-Host "Week 2 style touch-ups applied"

# ---------------------------------------------------------------------
# 3) Fill in the new Week-3 table (rows 40-57): header dates for all 7
#    days, but only the first 3 days (columns B:D) have data recorded.
# ---------------------------------------------------------------------
function Fill-RowPartial($rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, 2 + $i).Value = $values[$i]
    }
}

# Row 43: day headers (all 7 days get a header, even though only 3 have data)
Fill-RowPartial 43 @("第一天`n日期:2025-12-06","第二天`n日期:2025-12-07","第三天`n日期:2025-12-08","第四天`n日期:2025-12-9","第五天`n日期:2025-12-10","第六天`n日期:2025-12-11","第七天`n日期:2025-12-12")

Fill-RowPartial 44 @("7：30","7：40","7：30")
Fill-RowPartial 45 @("9：00","7：40","8：00")
Fill-RowPartial 46 @("23：00","23：10","23：00")
Fill-RowPartial 47 @("24：30","24：10","23：30")
Fill-RowPartial 48 @(20,30,10)
Fill-RowPartial 49 @(0,0,0)
Fill-RowPartial 50 @(0,0,0)
Fill-RowPartial 51 @(400,400,480)
Fill-RowPartial 52 @("无","无","无")
Fill-RowPartial 53 @(40,30,30)
Fill-RowPartial 54 @(4,4,3)
Fill-RowPartial 55 @(4,3,3)
Fill-RowPartial 56 @(4,4,4)
Fill-RowPartial 57 @("无","无","无")

Write-Host "Week 3 data filled"

# ---------------------------------------------------------------------
# 4) Update view state: scroll to the new block and leave the selection
#    on B56, and set up the page for (A4 / portrait) printing.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B56").Select()

$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait

Write-Host "View & page setup updated"
